# Add a new worksheet "ShopForRugs" at the end of the workbook (after
# "ShopForPaintings"), make it the active sheet, and populate it with the
# same A2/B2 "label / value" layout used by the other ShopFor* sheets.

$wb = $excel.ActiveWorkbook

$templateSheet = $wb.Worksheets.Item("ShopForPaintings")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# Add the new sheet right after the current last sheet so it lands at the
# end of the tab strip (and becomes the new active tab).
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "ShopForRugs"

# Reuse the existing cell formatting (fonts / number format / shared style
# indices) from the ShopForPaintings template instead of inventing new
# styles, by copying the formatted A2:B2 range first and then overwriting
# the values.
$templateSheet.Range("A2:B2").Copy()
$ws.Range("A2:B2").PasteSpecial(-4122)

$ws.Range("A2").Value = "4"
$ws.Range("B2").Value = "Traditional Area Rugs"

# Match the wider "value" column used on this sheet.
$ws.Columns.Item(2).ColumnWidth = 26.33203125

$ws.Range("B2").Select()
